# sudah membuat logic gugus
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "gugus" (group) values in column E for rows 2, 3, 6, 7 from 2 -> 1
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("E7").Value = 1

# Update the active selection to I9
$ws.Range("I9").Select()
